$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 7-15, column A: strip the " 00:00:00" time portion that was
# appended to these date strings, leaving just the date (as text).
# A leading apostrophe forces Excel to keep the value as text instead
# of re-parsing it into a date serial number.
$ws.Range("A7").Value  = "'2025-04-09"
$ws.Range("A8").Value  = "'2025-04-09"
$ws.Range("A9").Value  = "'2025-04-08"
$ws.Range("A10").Value = "'2025-04-08"
$ws.Range("A11").Value = "'2025-04-12"
$ws.Range("A12").Value = "'2025-04-14"
$ws.Range("A13").Value = "'2025-04-23"
$ws.Range("A14").Value = "'2025-04-23"
$ws.Range("A15").Value = "'2025-04-24"
